$wb = $excel.ActiveWorkbook

# 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7738
$ws1.Range("F4").Value = 7904
$ws1.Range("F8").Value = 6708
$ws1.Range("F9").Value = 6708
$ws1.Range("F10").Value = 3395
$ws1.Range("F12").Value = 3730
$ws1.Range("F14").Value = 52
$ws1.Range("F15").Value = 44
$ws1.Range("F16").Value = 68
$ws1.Range("F17").Value = 76
$ws1.Range("F18").Value = 473
$ws1.Range("F20").Value = 48
$ws1.Range("F21").Value = 323
$ws1.Range("F22").Value = 7
$ws1.Range("F23").Value = 333
$ws1.Range("F24").Value = 3877
$ws1.Range("F26").Value = 372
$ws1.Range("F28").Value = 286
$ws1.Range("F29").Value = 1496
$ws1.Range("F30").Value = 80
$ws1.Range("F33").Value = 1881
$ws1.Range("F35").Value = 50
$ws1.Range("F37").Value = 49
$ws1.Range("F38").Value = 3735
$ws1.Range("F39").Value = 331
$ws1.Range("F40").Value = 282
$ws1.Range("F41").Value = 44
$ws1.Range("F42").Value = 923
$ws1.Range("F43").Value = 545
$ws1.Range("F45").Value = 1438
$ws1.Range("F50").Value = 9

# 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 26
$ws2.Range("F6").Value = 415
$ws2.Range("F13").Value = 91
$ws2.Range("F17").Value = 87

# 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 137

# 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 137
$ws4.Range("F5").Value = 26
$ws4.Range("F7").Value = 7738
$ws4.Range("F9").Value = 7904
$ws4.Range("F12").Value = 6708
$ws4.Range("F13").Value = 3395
$ws4.Range("F15").Value = 3730
$ws4.Range("F17").Value = 52
$ws4.Range("F18").Value = 44
$ws4.Range("F19").Value = 68
$ws4.Range("F20").Value = 76
$ws4.Range("F21").Value = 473
$ws4.Range("F24").Value = 323
$ws4.Range("F25").Value = 333
$ws4.Range("F26").Value = 3877
$ws4.Range("F30").Value = 372
$ws4.Range("F32").Value = 1496
$ws4.Range("F33").Value = 80
$ws4.Range("F36").Value = 1881
$ws4.Range("F38").Value = 50
$ws4.Range("F40").Value = 91
$ws4.Range("F41").Value = 331
$ws4.Range("F42").Value = 282
$ws4.Range("F44").Value = 44
$ws4.Range("F45").Value = 923
$ws4.Range("F46").Value = 545
$ws4.Range("F47").Value = 87
